# Weekly update: insert a new "Sandia" price record as the new row 31
# (Vega Monumental Concepción, week of 2021-12-16), pushing the previously
# existing rows 31-64 down to 32-65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 31:64 down by one to make room for the new record.
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new weekly record.
$ws.Cells.Item(31, 1).Value = 11
$ws.Cells.Item(31, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(31, 3).Value = "Bíobío"
$ws.Cells.Item(31, 4).Value = 44546
$ws.Cells.Item(31, 5).Value = 8
$ws.Cells.Item(31, 6).Value = 100112028
$ws.Cells.Item(31, 7).Value = "Sandia"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 450
$ws.Cells.Item(31, 11).Value = 2500
$ws.Cells.Item(31, 12).Value = 3000
$ws.Cells.Item(31, 13).Value = 2778
$ws.Cells.Item(31, 14).Value = "$/unidad"
$ws.Cells.Item(31, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(31, 16).Value = 2778
$ws.Cells.Item(31, 17).Value = 1
$ws.Cells.Item(31, 18).Value = "Hortaliza"
